$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("testdata_Prop")
$ws2 = $wb.Worksheets.Item("testdata_Prop_g")

# ---------------------------------------------------------------------------
# Sheet "testdata_Prop": add rows 34-36 (Area9 / Area10 / Area11 with NAs)
# ---------------------------------------------------------------------------

# Row 34 - Area9, denominator only (numerator NA)
$ws1.Range("A34").Value = "Area9"
$ws1.Range("C34").Value = 100
$ws1.Range("G34").Value = "95%"
$ws1.Range("H34").Value = "proportion"
$ws1.Range("I34").Value = "Wilson"
$ws1.Range("J34").Value = 1

# Row 35 - Area10, numerator only (denominator NA)
$ws1.Range("A35").Value = "Area10"
$ws1.Range("B35").Value = 10
$ws1.Range("G35").Value = "95%"
$ws1.Range("H35").Value = "proportion"
$ws1.Range("I35").Value = "Wilson"
$ws1.Range("J35").Value = 1

# Row 36 - Area11, both numerator and denominator NA
$ws1.Range("A36").Value = "Area11"
$ws1.Range("G36").Value = "95%"
$ws1.Range("H36").Value = "proportion"
$ws1.Range("I36").Value = "Wilson"
$ws1.Range("J36").Value = 1

# ---------------------------------------------------------------------------
# Sheet "testdata_Prop_g": add rows 10-12 (grouped equivalents)
# ---------------------------------------------------------------------------

# Row 10 - Area9
$ws2.Range("A10").Value = "Area9"
$ws2.Range("C10").Value = 100
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = "95%"
$ws2.Range("G10").NumberFormat = "0%"
$ws2.Range("H10").Value = "proportion"
$ws2.Range("I10").Value = "Wilson"
$ws2.Range("J10").Value = 1

# Row 11 - Area10
$ws2.Range("A11").NumberFormat = "@"
$ws2.Range("A11").Value = "Area10"
$ws2.Range("A11").ClearFormats()
$ws2.Range("B11").Value = 10
$ws2.Range("G11").NumberFormat = "@"
$ws2.Range("G11").Value = "95%"
$ws2.Range("G11").ClearFormats()
$ws2.Range("H11").NumberFormat = "@"
$ws2.Range("H11").Value = "proportion"
$ws2.Range("H11").ClearFormats()
$ws2.Range("I11").NumberFormat = "@"
$ws2.Range("I11").Value = "Wilson"
$ws2.Range("I11").ClearFormats()
$ws2.Range("J11").Value = 1

# Row 12 - Area11
$ws2.Range("A12").NumberFormat = "@"
$ws2.Range("A12").Value = "Area11"
$ws2.Range("A12").ClearFormats()
$ws2.Range("G12").NumberFormat = "@"
$ws2.Range("G12").Value = "95%"
$ws2.Range("G12").ClearFormats()
$ws2.Range("H12").NumberFormat = "@"
$ws2.Range("H12").Value = "proportion"
$ws2.Range("H12").ClearFormats()
$ws2.Range("I12").NumberFormat = "@"
$ws2.Range("I12").Value = "Wilson"
$ws2.Range("I12").ClearFormats()
$ws2.Range("J12").Value = 1

# ---------------------------------------------------------------------------
# Sheet view / active tab: testdata_Prop becomes the active sheet
# ---------------------------------------------------------------------------

$ws2.Range("G27").Select() | Out-Null
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws1.Range("D38").Select() | Out-Null
